$d = $word.ActiveDocument

function Remove-ParagraphByExactText($exactText) {
    $count = $d.Paragraphs.Count
    for ($i = $count; $i -ge 1; $i--) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        # paragraph Range.Text includes the trailing paragraph mark char(s);
        # trim those for an exact compare against the visible text.
        $trimmed = $t.TrimEnd([char]13, [char]7)
        if ($trimmed -eq $exactText) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# 1) Remove the "nominal variable" and "ordinal variable" list paragraphs
#    (they followed "discrete variable").
Remove-ParagraphByExactText("nominal variable") | Out-Null
Remove-ParagraphByExactText("ordinal variable") | Out-Null

# 2) Remove the "biological vs technical replicates" list paragraph.
Remove-ParagraphByExactText("biological vs technical replicates") | Out-Null

# 3) Remove the "Correlation test" paragraph (its paragraph mark merges into
#    the following, otherwise-empty sectPr-bearing paragraph).
Remove-ParagraphByExactText("Correlation test") | Out-Null

# 4) Reword the fish reproductive-success example sentence.
$d.Content.Find.Execute(
    "You measure reproductive success of fish in your study you record sex and number of offspring that survive to adulthood for each fish in the study. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You measure the reproductive success of fish in your study and record the sex and number of offspring that survive to adulthood for each fish in the study. ",
    2) | Out-Null
